$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared string row for B13
$ws.Range("B13").Value = "Make recoil/accuracy dependent on whether or not player is standing still."

# Apply "Good" style to C4 and B5
$ws.Range("C4").Style = "Good"
$ws.Range("B5").Style = "Good"

# Apply "Neutral" style to B7
$ws.Range("B7").Style = "Neutral"

# Update selection to B14
$ws.Range("B14").Select()
